$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style of the existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for I and J columns (rows 2-9)
$values = @(
    @(7, 8),
    @(8, 8),
    @(6, 8),
    @(1, 2),
    @(8, 9),
    @(8, 9),
    @(8, 9),
    @(9, 9)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
